$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mistyped phone/document numbers on rows 8 and 9 (column F):
# 20999999999 -> 99999999999
$ws.Range("F8").Value = 99999999999
$ws.Range("F9").Value = 99999999999

# Update the sheet view: scroll so row 4 is at the top and select F10
# (mirrors the author re-opening/editing the receipt around that area).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F10").Select()
